$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").NumberFormat = "@"
$ws.Range("A2").Value = "José Pedro de Abulquerque"

$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "21999457635"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "jose@gmail.com"

$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "apartamento"

$ws.Range("F2").Value = 80
$ws.Range("G2").Value = 3

$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "50000"
